$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Book_01")

# Copy the formatting of the last existing data row (row 48) down across the
# new rows (49:59) so the new cells get the same style (s="2") as the rest
# of the table.
$ws.Range("A48:B48").Copy()
$ws.Range("A49:B59").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New Lane ID / Country (dirty) data appended to the table. Fill column A
# (Lane ID) fully first, then column B (Country), matching the order the
# shared strings table was built in the source edit.
$laneIds = @("A048", "A049", "A050", "A051", "A052", "A053", "A054", "A055", "A056", "A057", "A058")
$countries = @("Tokyo", "Aubervilliers", "Koper", "Kkorea", "Hawaï", "Norge", "Den russiske føderasjon", "Новая Зеландия", "pain", "Prtgal", "Zimbabyoue")

$row = 49
foreach ($laneId in $laneIds) {
    $ws.Cells.Item($row, 1).Value = $laneId
    $row = $row + 1
}

$row = 49
foreach ($country in $countries) {
    $ws.Cells.Item($row, 2).Value = $country
    $row = $row + 1
}

# Move the selection to the next empty cell below the newly-added rows,
# matching the selection left behind after entering the new data.
$ws.Range("B61").Select()
